$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

[void]$tr.InsertAfter("Infrastructure:")
[void]$tr.InsertAfter(" Multi-AZ cloud, auto-scaling, managed DB")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("Application:")
[void]$tr.InsertAfter(" Kubernetes, REST APIs, OAuth 2.0")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("Integration:")
[void]$tr.InsertAfter(" 5 enterprise systems connected")

$labels = @("Infrastructure:", "Application:", "Integration:")
for ($i = 1; $i -le 3; $i++) {
    $para = $tr.Paragraphs($i)
    $bold = $tr.Characters($para.Start, $labels[$i - 1].Length)
    $bold.Font.Bold = $true
}
